# Update functions for VRMS/IRMS register and Energy register, and add
# new Alert / threshold registers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 26: "V and I" register (Reg 14) ---
# Length(Bytes) 18 -> 12, description gets an arrow added ("-> AC and BC")
$ws.Range("C26").Value = 12
$ws.Range("E26").Value = "3*2*2 Bytes (VRMS, IRMS ) (3 bytes each, two phase each -> AC and BC)"

# --- Row 27: "Energy" register (Reg 15) ---
# Length(Bytes) 18 -> 12, description gets an arrow added ("-> AC and BC")
$ws.Range("C27").Value = 12
$ws.Range("E27").Value = "2*3*2 Bytes (Whr,VARhr,Vahr) (2 bytes each, two phase each -> AC and BC)"

# --- Row 28 (new): Alert byte register ---
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = "R"
$ws.Range("E28").Value = "1 byte Alert[bit7-0]=(overVoltage AC, overVoltage BC, under Voltage AC, under Voltage BC, overCurrent A, overCurrent B, 0, 0)"

# --- Row 29: append threshold columns next to the existing "<! Regular Registers>" marker ---
$ws.Range("C29").Value = 18
$ws.Range("D29").Value = "R/W"
$ws.Range("E29").Value = "threshold of 3 bytes each corresponding to the above alerts."

# --- Update the saved view state (scroll position / selection) ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("E32").Select()
